$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "73.60") and Excel would auto-convert them
# to real numbers on a General-formatted cell, dropping trailing zeros / using the
# native decimal representation. The source data stores these as literal text, so
# force Text format before writing, then clear the format change back off so the
# cell keeps its original (default) style - matching the source workbook exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.529.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.487.00'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.63%  '

$ws.Range("E7").Value = '  -1.00%  '

$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.62'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.34%  '

$ws.Range("E11").Value = '  +0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.867.91'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("E14").Value = '  -2.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.79'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.476.92'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.753'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.551.13'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.35'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0929'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.24'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.25'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.08'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.44%  '

$ws.Range("E24").Value = '  -2.97%  '

$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  -1.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.96'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.85%  '

$ws.Range("E28").Value = '  -0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.27'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '157.31'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.44'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.93%  '

$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.13'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.03%  '

$ws.Range("E35").Value = '  -0.31%  '

$ws.Range("E36").Value = '  -5.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.93'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.106'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.48%  '

$ws.Range("E39").Value = '  -3.00%  '

$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.15'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.53%  '

$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.87'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -6.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.964.88'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.98'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.86'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.726.92'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.52'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.59'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.60'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.55%  '

